$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = "Climate Action, SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    30 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    34 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    35 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    36 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    38 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    40 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    41 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    44 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    47 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    48 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    49 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    51 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    53 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    55 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    56 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    57 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    58 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    59 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    61 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    62 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    63 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    64 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    65 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    80 = "SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    87 = "Inclusive Education, SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    88 = "Agriculture, SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    89 = "Agriculture, SDG 10 (Reduced inequality through access to info in local languages/NLP)"
    90 = "Agriculture, SDG 10 (Reduced inequality through access to info in local languages/NLP)"
}

foreach ($row in $newValues.Keys) {
    $ws.Range("H$row").Value = $newValues[$row]
}
